# Updates the cryptos list (Price / Volume(1h) columns) per the scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "48.092.23"; E = "+0.16%" },
    @{ Row = 3; D = "2.498.56"; E = "-0.49%" },
    @{ Row = 4; D = $null; E = "-0.01%" },
    @{ Row = 5; D = "320.78"; E = "-0.58%" },
    @{ Row = 6; D = "107.53"; E = "-1.87%" },
    @{ Row = 7; D = $null; E = "+0.13%" },
    @{ Row = 8; D = "0.999"; E = "-0.08%" },
    @{ Row = 9; D = "0.539"; E = "-1.69%" },
    @{ Row = 10; D = "39.61"; E = "-2.18%" },
    @{ Row = 11; D = "20.19"; E = "+8.06%" },
    @{ Row = 12; D = "0.0813"; E = "-0.03%" },
    @{ Row = 13; D = $null; E = "-0.30%" },
    @{ Row = 14; D = "7.12"; E = "-2.05%" },
    @{ Row = 15; D = "2.889.65"; E = "-0.31%" },
    @{ Row = 16; D = "2.501.87"; E = "-0.84%" },
    @{ Row = 17; D = "0.837"; E = "-1.92%" },
    @{ Row = 18; D = "47.937.87"; E = "+0.10%" },
    @{ Row = 19; D = "12.99"; E = "-2.56%" },
    @{ Row = 20; D = "6.73"; E = "+0.98%" },
    @{ Row = 21; D = "0.0₃0940"; E = "-0.63%" },
    @{ Row = 22; D = "2.74"; E = "-1.24%" },
    @{ Row = 23; D = "277.59"; E = "+11.87%" },
    @{ Row = 24; D = "71.73"; E = "+1.24%" },
    @{ Row = 25; D = "2.54"; E = "-1.01%" },
    @{ Row = 26; D = $null; E = "-0.12%" },
    @{ Row = 27; D = "25.65"; E = "-1.16%" },
    @{ Row = 28; D = "9.75"; E = "-3.03%" },
    @{ Row = 29; D = $null; E = "+0.27%" },
    @{ Row = 30; D = "35.01"; E = "-0.32%" },
    @{ Row = 31; D = $null; E = "-4.56%" },
    @{ Row = 32; D = "49.41"; E = "-0.80%" },
    @{ Row = 33; D = "19.61"; E = "-3.10%" },
    @{ Row = 34; D = "1.00"; E = "-0.14%" },
    @{ Row = 35; D = "5.31"; E = "-1.37%" },
    @{ Row = 36; D = "0.0779"; E = "-1.38%" },
    @{ Row = 37; D = $null; E = "-1.85%" },
    @{ Row = 38; D = "4.63"; E = "-1.90%" },
    @{ Row = 39; D = "2.89"; E = "-3.04%" },
    @{ Row = 40; D = $null; E = "-0.80%" },
    @{ Row = 41; D = "120.64"; E = "+1.03%" },
    @{ Row = 42; D = "2.21"; E = "-0.35%" },
    @{ Row = 43; D = "21.38"; E = "-5.48%" },
    @{ Row = 44; D = "0.0301"; E = "+0.56%" },
    @{ Row = 45; D = "2.010.04"; E = "+0.43%" },
    @{ Row = 46; D = "3.15"; E = "+2.75%" },
    @{ Row = 47; D = $null; E = "-2.14%" },
    @{ Row = 48; D = "1.85"; E = "-1.44%" },
    @{ Row = 49; D = "8.97"; E = "-1.17%" },
    @{ Row = 50; D = "5.16"; E = "-1.62%" },
    @{ Row = 51; D = "80.15"; E = "+2.58%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the value as text even
        # when it looks numeric, then reset the style so no quote-prefix
        # formatting marker lingers on the cell.
        $ws.Range("D$r").Value = "'" + $u.D
        $ws.Range("D$r").Style = "Normal"
    }
    $ws.Range("E$r").Value = "  " + $u.E + "  "
}
